$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting the existing rows 45-100 down to 46-101.
$ws.Rows("45:45").Insert()

# Populate the newly inserted row 45 with the new weekly record.
$ws.Cells.Item(45, 1).Value  = 10
$ws.Cells.Item(45, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(45, 3).Value  = "La Araucanía"
$ws.Cells.Item(45, 4).Value  = 44413
$ws.Cells.Item(45, 5).Value  = 9
$ws.Cells.Item(45, 6).Value  = 100112005
$ws.Cells.Item(45, 7).Value  = "Puerro"
$ws.Cells.Item(45, 8).Value  = "Azul de Maquehue"
$ws.Cells.Item(45, 9).Value  = "Primera"
$ws.Cells.Item(45, 10).Value = 50
$ws.Cells.Item(45, 11).Value = 8000
$ws.Cells.Item(45, 12).Value = 8000
$ws.Cells.Item(45, 13).Value = 8000
$ws.Cells.Item(45, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(45, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(45, 16).Value = 667
$ws.Cells.Item(45, 17).Value = 12
$ws.Cells.Item(45, 18).Value = "Hortaliza"
